$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block: card holder name ---
$ws.Range("C2").Value = "Hartmut"

# --- Card number: must stay text, not get auto-coerced to a number by the
#     "looks like a number" heuristic that a plain .Value assignment
#     triggers. Build it as a text formula, then paste-special as values
#     only so the literal text lands back in the cell without dragging a
#     new "Text" number format style along with it (keeps the original
#     cell style index intact, same as a manual Excel copy/paste-values
#     workflow). ---
$ws.Range("B3").Formula = "=""2570314725427075"""
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 21.04.2024"

# --- Row 6: MITGLIEDSBEITRAG ZEUS BODYPOWER transaction (dates + amount only change) ---
$ws.Range("B6").Value = "25.04."
$ws.Range("C6").Value = "26.04."
$ws.Range("E6").Value = "24,83-"

# --- Row 7: transaction replaced ---
$ws.Range("B7").Value = "26.04."
$ws.Range("C7").Value = "27.04."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU XVBLQL"
$ws.Range("E7").Value = "19,61-"

# --- Row 8: transaction replaced ---
$ws.Range("B8").Value = "30.04."
$ws.Range("C8").Value = "01.05."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-44726137"
$ws.Range("E8").Value = "55,71-"

# --- Row 9: transaction removed, cells cleared back to blank
#     (E9's alignment switches to centered to match the blank-row style).
#     NB: assigning "" clears reliably even on the merged E9:F9 cell;
#     .ClearContents() was observed to silently no-op on that merged range. ---
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# --- Row 10: transaction removed, cells cleared back to blank
#     (E10's alignment switches to right/center to match the blank-row style) ---
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 05.05.2024"
$ws.Range("E12").Value = "100,15-"

# --- Next billing date note ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 14.05.2024"
